$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.185544452515558
$ws.Range("C2").Value2 = 0.06383851775807869
$ws.Range("D2").Value2 = 0.07727191630942798
$ws.Range("E2").Value2 = 0.08551212623796545
$ws.Range("G2").Value2 = 1.865936803006377
$ws.Range("H2").Value2 = 1.591832250881993
$ws.Range("K2").Value2 = 0.7384429030062734
$ws.Range("L2").Value2 = 0.2206154414127894
$ws.Range("N2").Value2 = 2.710617785907999

$ws.Range("B3").Value2 = 1.137594777445287
$ws.Range("C3").Value2 = 0.06071183784505507
$ws.Range("D3").Value2 = 0.07027250815387731
$ws.Range("E3").Value2 = 0.08496741812924924
$ws.Range("G3").Value2 = 1.850413626096724
$ws.Range("H3").Value2 = 1.590036896018105
$ws.Range("K3").Value2 = 0.6928460921953103
$ws.Range("L3").Value2 = 0.2136260789454951
$ws.Range("N3").Value2 = 2.723139883655065

$ws.Range("B4").Value2 = 1.108845425852991
$ws.Range("C4").Value2 = 0.05876537745957933
$ws.Range("D4").Value2 = 0.0660127414868299
$ws.Range("E4").Value2 = 0.08467119388705768
$ws.Range("G4").Value2 = 1.841801578447374
$ws.Range("H4").Value2 = 1.589565559122917
$ws.Range("K4").Value2 = 0.6652736217377253
$ws.Range("L4").Value2 = 0.2094638957522079
$ws.Range("N4").Value2 = 2.731549818837543

$ws.Range("B5").Value2 = 1.097303835280002
$ws.Range("C5").Value2 = 0.05796543455799252
$ws.Range("D5").Value2 = 0.06428631605588464
$ws.Range("E5").Value2 = 0.08456009926772623
$ws.Range("G5").Value2 = 1.838522883859213
$ws.Range("H5").Value2 = 1.589532124424323
$ws.Range("K5").Value2 = 0.654144161821705
$ws.Range("L5").Value2 = 0.207800262874386
$ws.Range("N5").Value2 = 2.735158143586794

$ws.Range("B6").Value2 = 1.095397872832507
$ws.Range("C6").Value2 = 0.0578321961614634
$ws.Range("D6").Value2 = 0.06400021412764545
$ws.Range("E6").Value2 = 0.08454223325322907
$ws.Range("G6").Value2 = 1.837992387339952
$ws.Range("H6").Value2 = 1.589536153050787
$ws.Range("K6").Value2 = 0.6523025554862443
$ws.Range("L6").Value2 = 0.2075259795126385
$ws.Range("N6").Value2 = 2.735768242670858

$ws.Range("B7").Value2 = 1.108689067377412
$ws.Range("C7").Value2 = 0.05875461650934
$ws.Range("D7").Value2 = 0.0659894200755673
$ws.Range("E7").Value2 = 0.0846696566698828
$ws.Range("G7").Value2 = 1.841756426878277
$ws.Range("H7").Value2 = 1.589564465943425
$ws.Range("K7").Value2 = 0.665123094693854
$ws.Range("L7").Value2 = 0.2094413278975225
$ws.Range("N7").Value2 = 2.731597748506204

$ws.Range("B8").Value2 = 1.168867777269497
$ws.Range("C8").Value2 = 0.06276594883755138
$ws.Range("D8").Value2 = 0.07485061274189775
$ws.Range("E8").Value2 = 0.08531637983508489
$ws.Range("G8").Value2 = 1.860393390955295
$ws.Range("H8").Value2 = 1.591082218412907
$ws.Range("K8").Value2 = 0.7226330046054272
$ws.Range("L8").Value2 = 0.2181786513121295
$ws.Range("N8").Value2 = 2.714785590945525

$ws.Range("B9").Value2 = 1.292375943890704
$ws.Range("C9").Value2 = 0.07042263748118671
$ws.Range("D9").Value2 = 0.09253199035607906
$ws.Range("E9").Value2 = 0.086887825260618
$ws.Range("G9").Value2 = 1.904254968994053
$ws.Range("H9").Value2 = 1.599068918329777
$ws.Range("K9").Value2 = 0.8387872180522891
$ws.Range("L9").Value2 = 0.236340971679752
$ws.Range("N9").Value2 = 2.687548424322003

$ws.Range("B10").Value2 = 1.386491262260108
$ws.Range("C10").Value2 = 0.07592356714809512
$ws.Range("D10").Value2 = 0.1057149400281077
$ws.Range("E10").Value2 = 0.08822733005385075
$ws.Range("G10").Value2 = 1.940973984157552
$ws.Range("H10").Value2 = 1.607998921978862
$ws.Range("K10").Value2 = 0.9262112725014049
$ws.Range("L10").Value2 = 0.2503168996042859
$ws.Range("N10").Value2 = 2.671042956265367

$ws.Range("B11").Value2 = 1.430044947341173
$ws.Range("C11").Value2 = 0.07839979564995758
$ws.Range("D11").Value2 = 0.1117555270147932
$ws.Range("E11").Value2 = 0.08887691188329327
$ws.Range("G11").Value2 = 1.958662173005479
$ws.Range("H11").Value2 = 1.612728236764923
$ws.Range("K11").Value2 = 0.9664419930385293
$ws.Range("L11").Value2 = 0.2568133777561741
$ws.Range("N11").Value2 = 2.664298043872776

$ws.Range("B12").Value2 = 1.446644245320101
$ws.Range("C12").Value2 = 0.07933376459101282
$ws.Range("D12").Value2 = 0.1140493038042933
$ws.Range("E12").Value2 = 0.08912867618441567
$ws.Range("G12").Value2 = 1.965502337936044
$ws.Range("H12").Value2 = 1.614615132418436
$ws.Range("K12").Value2 = 0.9817429332529741
$ws.Range("L12").Value2 = 0.2592934413559931
$ws.Range("N12").Value2 = 2.661853940706564

$ws.Range("B13").Value2 = 1.443064554007833
$ws.Range("C13").Value2 = 0.07913278272702939
$ws.Range("D13").Value2 = 0.1135550149730591
$ws.Range("E13").Value2 = 0.08907419715251663
$ws.Range("G13").Value2 = 1.964022861774765
$ws.Range("H13").Value2 = 1.614204484254884
$ws.Range("K13").Value2 = 0.978444644861554
$ws.Range("L13").Value2 = 0.2587584257325091
$ws.Range("N13").Value2 = 2.662375423456112

$ws.Range("B14").Value2 = 1.431408447714261
$ws.Range("C14").Value2 = 0.0784767082712392
$ws.Range("D14").Value2 = 0.1119441100038756
$ws.Range("E14").Value2 = 0.08889750883048464
$ws.Range("G14").Value2 = 1.959222068319775
$ws.Range("H14").Value2 = 1.612881548151165
$ws.Range("K14").Value2 = 0.9676994782506938
$ws.Range("L14").Value2 = 0.2570170132815122
$ws.Range("N14").Value2 = 2.66409475900403

$ws.Range("B15").Value2 = 1.424282611610465
$ws.Range("C15").Value2 = 0.07807435975620081
$ws.Range("D15").Value2 = 0.110958211575138
$ws.Range("E15").Value2 = 0.08879003497484561
$ws.Range("G15").Value2 = 1.956299954123494
$ws.Range("H15").Value2 = 1.61208371727767
$ws.Range("K15").Value2 = 0.9611264099698076
$ws.Range("L15").Value2 = 0.255952951997088
$ws.Range("N15").Value2 = 2.665162240981672

$ws.Range("B16").Value2 = 1.383659822677487
$ws.Range("C16").Value2 = 0.07576121642719613
$ws.Range("D16").Value2 = 0.1053210565852396
$ws.Range("E16").Value2 = 0.08818568763528134
$ws.Range("G16").Value2 = 1.939837869615161
$ws.Range("H16").Value2 = 1.607703279836073
$ws.Range("K16").Value2 = 0.9235913899755985
$ws.Range("L16").Value2 = 0.2498951346262288
$ws.Range("N16").Value2 = 2.671499137770255

$ws.Range("B17").Value2 = 1.358928647656398
$ws.Range("C17").Value2 = 0.07433550193979954
$ws.Range("D17").Value2 = 0.1018740593880381
$ws.Range("E17").Value2 = 0.08782524228022126
$ws.Range("G17").Value2 = 1.929991398664441
$ws.Range("H17").Value2 = 1.605186924837625
$ws.Range("K17").Value2 = 0.9006830602115201
$ws.Range("L17").Value2 = 0.246214423871308
$ws.Range("N17").Value2 = 2.675582353016011

$ws.Range("B18").Value2 = 1.344773594609308
$ws.Range("C18").Value2 = 0.07351300293655072
$ws.Range("D18").Value2 = 0.0998955372266721
$ws.Range("E18").Value2 = 0.08762171091863635
$ws.Range("G18").Value2 = 1.924420609358179
$ws.Range("H18").Value2 = 1.603802360541096
$ws.Range("K18").Value2 = 0.8875501853418371
$ws.Range("L18").Value2 = 0.2441104348514358
$ws.Range("N18").Value2 = 2.678002759679217

$ws.Range("B19").Value2 = 1.339992901057769
$ws.Range("C19").Value2 = 0.07323409424290617
$ws.Range("D19").Value2 = 0.09922634562441601
$ws.Range("E19").Value2 = 0.08755344928589182
$ws.Range("G19").Value2 = 1.922550333027004
$ws.Range("H19").Value2 = 1.603344350078601
$ws.Range("K19").Value2 = 0.8831110678433163
$ws.Range("L19").Value2 = 0.243400302070313
$ws.Range("N19").Value2 = 2.678834600217812

$ws.Range("B20").Value2 = 1.361554113229658
$ws.Range("C20").Value2 = 0.07448752659914248
$ws.Range("D20").Value2 = 0.1022405733530434
$ws.Range("E20").Value2 = 0.08786322036414518
$ws.Range("G20").Value2 = 1.931029981068264
$ws.Range("H20").Value2 = 1.605448297326717
$ws.Range("K20").Value2 = 0.9031171992955933
$ws.Range("L20").Value2 = 0.2466048900639635
$ws.Range("N20").Value2 = 2.675140249790189

$ws.Range("B21").Value2 = 1.434829239048838
$ws.Range("C21").Value2 = 0.07866951403133271
$ws.Range("D21").Value2 = 0.1124170992957971
$ws.Range("E21").Value2 = 0.08894924958860173
$ws.Range("G21").Value2 = 1.96062831935069
$ws.Range("H21").Value2 = 1.613267520209945
$ws.Range("K21").Value2 = 0.9708537898871725
$ws.Range("L21").Value2 = 0.2575279655703184
$ws.Range("N21").Value2 = 2.663586759710327

$ws.Range("B22").Value2 = 1.483339333864649
$ws.Range("C22").Value2 = 0.08138099010430722
$ws.Range("D22").Value2 = 0.1191050196368337
$ws.Range("E22").Value2 = 0.08969273130942312
$ws.Range("G22").Value2 = 1.980800630570911
$ws.Range("H22").Value2 = 1.618937480377866
$ws.Range("K22").Value2 = 1.015510894838485
$ws.Range("L22").Value2 = 0.264783350976586
$ws.Range("N22").Value2 = 2.656677457164761

$ws.Range("B23").Value2 = 1.457391790257191
$ws.Range("C23").Value2 = 0.07993579658405281
$ws.Range("D23").Value2 = 0.115532145015834
$ws.Range("E23").Value2 = 0.08929283891667339
$ws.Range("G23").Value2 = 1.969958365624109
$ws.Range("H23").Value2 = 1.615860077165536
$ws.Range("K23").Value2 = 0.9916410704493899
$ws.Range("L23").Value2 = 0.2609003414564484
$ws.Range("N23").Value2 = 2.66030629200344

$ws.Range("B24").Value2 = 1.360366943294537
$ws.Range("C24").Value2 = 0.07441880509099974
$ws.Range("D24").Value2 = 0.10207486242291
$ws.Range("E24").Value2 = 0.08784603896784304
$ws.Range("G24").Value2 = 1.930560157399128
$ws.Range("H24").Value2 = 1.605329937309165
$ws.Range("K24").Value2 = 0.9020166084420111
$ws.Range("L24").Value2 = 0.2464283226116208
$ws.Range("N24").Value2 = 2.675339897534442

$ws.Range("B25").Value2 = 1.258372853293281
$ws.Range("C25").Value2 = 0.06837339157451083
$ws.Range("D25").Value2 = 0.08771540307448333
$ws.Range("E25").Value2 = 0.08643023469967659
$ws.Range("G25").Value2 = 1.891602603748026
$ws.Range("H25").Value2 = 1.596371124214386
$ws.Range("K25").Value2 = 0.8070001188514766
$ws.Range("L25").Value2 = 0.2313169560587909
$ws.Range("N25").Value2 = 2.694302024032382
